# Applies the "working on case study edits" commit to slide 1 of the
# Magallanes case-study deck.
#
# Summary of changes (see xml_diff):
#   1. "Rounded Rectangle 4" (id 5, the lower big outline box) shifts
#      up/left a bit (offset only, size unchanged).
#   2. "TextBox 77" (id 78, "Shared Data sets") grows taller and moves up
#      slightly (no text change).
#   3. "TextBox 78" (id 79, "Scrapped Data sets") grows taller/moves up
#      AND its text typo is corrected to "Scraped Data sets" (PowerPoint
#      records the fix as a new leading run "Scraped " followed by the
#      untouched "Data sets" run).
#   4. "TextBox 80" (id 81, "PYTHON Program") grows taller/moves up AND
#      its text changes to "Python Program".
#   5 & 6. The two connector arrows around boxes 79/81 (id 82 and id 85)
#      are re-routed (new Top/Height) to match the resized boxes.
#   7. "TextBox 127" (id 128, "GitHub Desktop") merges its " " and
#      "Desktop" runs into a single " Desktop" run.
#
# NOTE on geometry: Shape.Left/Top/Width/Height are expressed in points
# (1 pt = 12700 EMU) and are stored as single-precision floats internally,
# so a plain "EMU / 12700.0" literal can truncate one EMU short once it is
# converted back to EMU on save. The literals below were picked so they
# round-trip to the exact target EMU values from the diff.

function Get-ShapeById {
    param($Shapes, [int]$Id)
    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $candidate = $Shapes.Item($i)
        if ($candidate.Id -eq $Id) {
            return $candidate
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) Rounded Rectangle 4 (id 5) - offset only, size unchanged.
$sh = Get-ShapeById $s.Shapes 5
$sh.Left = 56.34252166748047
$sh.Top = 175.83685302734375

# 2) TextBox 77 (id 78) "Shared Data sets" - resize only, text unchanged.
$sh = Get-ShapeById $s.Shapes 78
$sh.Top = 346.1612854003906
$sh.Height = 37.53748321533203

# 3) TextBox 78 (id 79) "Scrapped Data sets" -> "Scraped Data sets"
#    (resize + fix typo, which splits the run into "Scraped " + "Data sets").
$sh = Get-ShapeById $s.Shapes 79
$sh.Top = 346.2223815917969
$sh.Height = 37.53748321533203
$sub = $sh.TextFrame.TextRange.Characters(1, 9)
$sub.Text = "Scraped "

# 4) TextBox 80 (id 81) "PYTHON Program" -> "Python Program" (resize + text).
$sh = Get-ShapeById $s.Shapes 81
$sh.Top = 403.0401611328125
$sh.Height = 37.53748321533203
$sh.TextFrame.TextRange.Text = "Python Program"

# 5) Straight Arrow Connector 81 (id 82) - resize to follow box 81's move.
$sh = Get-ShapeById $s.Shapes 82
$sh.Top = 440.5776672363281
$sh.Height = 23.648977279663086

# 6) Straight Arrow Connector 84 (id 85) - resize to follow box 79's move.
$sh = Get-ShapeById $s.Shapes 85
$sh.Top = 383.7598571777344
$sh.Height = 19.280315399169922

# 7) TextBox 127 (id 128) "GitHub" + " " + "Desktop" -> "GitHub" + " Desktop".
$sh = Get-ShapeById $s.Shapes 128
$sub = $sh.TextFrame.TextRange.Characters(7, 8)
$sub.Text = " Desktop"
